# Weekly refresh: a new week of data (Fecha=44651) is inserted at the top of
# the data block (rows 156-157), pushing every existing week down by one
# slot (two rows). The last existing week (rows 254-255, Fecha=44628) ends
# up re-homed at the new bottom rows 256-257, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the existing data block (rows 156-255) before it gets
#    overwritten, then re-write it two rows lower (rows 158-257).
$srcRange = $ws.Range("A156:R255")
$block = $srcRange.Value2

$dstRange = $ws.Range("A158:R257")
$dstRange.Value2 = $block

# 2) Also copy row formatting (the date column D carries a date style) down
#    two rows so the newly populated rows 158-257 keep the same look as
#    their source rows did.
$srcRange.Copy()
$dstRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Populate the two brand-new rows (156-157) for the new week
#    (Fecha 44651), keeping everything else about those records the same
#    shape as the rows that used to sit there, and overwriting only the
#    fields that actually changed.
$ws.Range("A156").Value2 = 11
$ws.Range("B156").Value2 = "Vega Monumental Concepción"
$ws.Range("C156").Value2 = "Bíobío"
$ws.Range("D156").Value2 = 44651
$ws.Range("E156").Value2 = 8
$ws.Range("F156").Value2 = 100112023
$ws.Range("G156").Value2 = "Brócoli"
$ws.Range("H156").Value2 = "Sin especificar"
$ws.Range("I156").Value2 = "Primera"
$ws.Range("J156").Value2 = 2000
$ws.Range("K156").Value2 = 900
$ws.Range("L156").Value2 = 1000
$ws.Range("M156").Value2 = 950
$ws.Range("N156").Value2 = "`$/unidad"
$ws.Range("O156").Value2 = "Región Metropolitana"
$ws.Range("P156").Value2 = 950
$ws.Range("Q156").Value2 = 1
$ws.Range("R156").Value2 = "Hortaliza"

$ws.Range("A157").Value2 = 11
$ws.Range("B157").Value2 = "Vega Monumental Concepción"
$ws.Range("C157").Value2 = "Bíobío"
$ws.Range("D157").Value2 = 44651
$ws.Range("E157").Value2 = 8
$ws.Range("F157").Value2 = 100112023
$ws.Range("G157").Value2 = "Brócoli"
$ws.Range("H157").Value2 = "Sin especificar"
$ws.Range("I157").Value2 = "Segunda"
$ws.Range("J157").Value2 = 1500
$ws.Range("K157").Value2 = 800
$ws.Range("L157").Value2 = 800
$ws.Range("M157").Value2 = 800
$ws.Range("N157").Value2 = "`$/unidad"
$ws.Range("O157").Value2 = "Región Metropolitana"
$ws.Range("P157").Value2 = 800
$ws.Range("Q157").Value2 = 1
$ws.Range("R157").Value2 = "Hortaliza"
